# backwardElimination.xlsx holds 16 worksheets (tabs "23" down to "8"), one
# per step of a backward-elimination MLR run. Each sheet's B2 cell holds a
# single text blob: the verbatim statsmodels "OLS Regression Results"
# console printout for that step, including a "Date:" line and a "Time:"
# line stamped when the Python script that produced the workbook was run.
#
# The workbook was regenerated by re-running that script on a later date
# (Sat, 28 Dec 2019 -> Sun, 29 Dec 2019) and at a different time-of-day
# (20:59:45 -> 16:11:14). Every other figure in every summary block is
# unchanged (same R-squared, coefficients, AIC/BIC, etc.), so the edit is a
# pure find/replace of those two stamps inside every worksheet's B2 cell.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    # NOTE: use Value2 (Value's getter on this host is a stub); Value2 also
    # round-trips the full multi-line text faithfully on write.
    $text = $cell.Value2
    if ($null -ne $text) {
        $newText = $text -replace "Sat, 28 Dec 2019", "Sun, 29 Dec 2019"
        $newText = $newText -replace "20:59:45", "16:11:14"
        if ($newText -ne $text) {
            $cell.Value2 = $newText
        }
    }
}
